$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "XII ..." (12th grade) class rows that are no longer needed.
# This clears the id/kode/kelas data in rows 20-27 while leaving the
# existing cell formatting (e.g. style on column C) intact.
$ws.Range("A20:C27").ClearContents()

# Update the active selection to reflect where the user left off editing.
$ws.Range("C12").Select()
